# Update "想去人数" (number of people interested) counts in the
# "展览" and "全部类型" worksheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 2,3,5,6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1487
$ws1.Range("F3").Value = 3133
$ws1.Range("F5").Value = 843
$ws1.Range("F6").Value = 295

# Sheet "全部类型" (All types) - rows 2,3,5,7
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value = 1487
$ws2.Range("F3").Value = 3133
$ws2.Range("F5").Value = 843
$ws2.Range("F7").Value = 295
